$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/E (and non-float-like D values) can be set directly as text.
# Float-looking D values (e.g. "249.67") need NumberFormat forced to text first,
# otherwise Excel auto-converts them to numbers, losing the intended text content.

# Row 2
$ws.Range("D2").Value = '37.108.25'
$ws.Range("E2").Value = '  +0.30%  '

# Row 3
$ws.Range("D3").Value = '2.054.45'
$ws.Range("E3").Value = '  +0.04%  '

# Row 4
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.67'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +0.47%  '

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.667'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.20'
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = '  +10.80%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  +1.94%  '

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0791'
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  +0.01%  '

# Row 11
$ws.Range("E11").Value = '  +2.10%  '

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.16'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  +8.28%  '

# Row 13
$ws.Range("D13").Value = '2.355.97'
$ws.Range("E13").Value = '  +0.12%  '

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.821'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  +1.26%  '

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.73'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  +10.02%  '

# Row 16
$ws.Range("D16").Value = '2.057.34'
$ws.Range("E16").Value = '  +0.21%  '

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.60'
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  +31.50%  '

# Row 18
$ws.Range("D18").Value = '37.073.09'
$ws.Range("E18").Value = '  +0.39%  '

# Row 19
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '75.40'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  +4.24%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0905'
$ws.Range("E20").Value = '  -4.19%  '

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.43'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  +2.44%  '

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.01'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  +0.86%  '

# Row 23
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -0.25%  '

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  +12.07%  '

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.90'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  -0.52%  '

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.38'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  +4.56%  '

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.99'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -0.02%  '

# Row 29
$ws.Range("E29").Value = '  +1.21%  '

# Row 30
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.15'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  +11.35%  '

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.85'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  +7.21%  '

# Row 32
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0621'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  +0.31%  '

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.52'
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  +4.86%  '

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0887'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  +3.03%  '

# Row 35
$ws.Range("E35").Value = '  +0.05%  '

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.26'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  +0.12%  '

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.74'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -1.32%  '

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.109'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  +4.28%  '

# Row 39
$ws.Range("E39").Value = '  +0.84%  '

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.31'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  +31.06%  '

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.15'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  +14.28%  '

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '18.00'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  +0.74%  '

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0224'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  +0.76%  '

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.14'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  +0.47%  '

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '97.60'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  +1.81%  '

# Row 46
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.46'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +3.43%  '

# Row 47
$ws.Range("D47").Value = '1.293.20'
$ws.Range("E47").Value = '  +0.21%  '

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.84'
$ws.Range("D48").Style = $style
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("E48").Value = '  -9.17%  '

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.88'
$ws.Range("D49").Style = $style
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("E49").Value = '  -1.63%  '

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.86'
$ws.Range("D50").Style = $style
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E50").Value = '  +1.15%  '

# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.242.33'
$ws.Range("E51").Value = '  +0.22%  '
